$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value  = "Financial conditions & employment, US.xlsx"
$ws.Range("A3").Value  = "Bitty, SPX & GM2 fitted trends.xlsx"
$ws.Range("A4").Value  = "Monetary Expansion and Risk Asset Activity.xlsx"
$ws.Range("A5").Value  = "Personal savings U.S (BEA).xlsx"
$ws.Range("A6").Value  = "Bitcoin price change is fueled by global monetary growth.xlsx"
$ws.Range("A7").Value  = "U.S GDP and GDI.xlsx"
$ws.Range("A8").Value  = "Real Incomes U.S (BEA).xlsx"
$ws.Range("A9").Value  = "Bitcoin price change is fueled by global monetary growth (with forecast).xlsx"
$ws.Range("A10").Value = "US financial conditions, employment & equities.xlsx"
$ws.Range("A11").Value = "US Equity Indices.xlsx"
$ws.Range("A12").Value = "US Gross Domestic Income and Equity Indices.xlsx"
$ws.Range("A13").Value = "US and Global Monetary Aggregates.xlsx"
$ws.Range("A14").Value = "Bank credit and M2 U.S.xlsx"
$ws.Range("A16").Value = "Monetary Aggregates & Inflation, USA.xlsx"
